# Fruta / hortaliza, semanal
#
# A new weekly price record for Naranja (Lane Late, Primera) at
# Terminal Hortofrutícola Agro Chillán needs to be inserted ahead of the
# existing row 438, pushing the existing rows 438-468 down to 439-469.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 438; everything that used
# to live at row 438 (and below, down to 468) shifts down to 439 (..469).
$ws.Rows(438).Insert()

# Populate the newly inserted row 438 with the new record.
$ws.Range("A438").Value = 7
$ws.Range("B438").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C438").Value = "Ñuble"
$ws.Range("D438").Value = 44826
$ws.Range("E438").Value = 16
$ws.Range("F438").Value = "Fruta"
$ws.Range("G438").Value = 100102
$ws.Range("H438").Value = "Cítricos"
$ws.Range("I438").Value = 100102005
$ws.Range("J438").Value = "Naranja"
$ws.Range("K438").Value = "Lane Late"
$ws.Range("L438").Value = "Primera"
$ws.Range("M438").Value = 160
$ws.Range("N438").Value = 6500
$ws.Range("O438").Value = 7000
$ws.Range("P438").Value = 6750
$ws.Range("Q438").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R438").Value = "Región de O'Higgins"
$ws.Range("S438").Value = 450
$ws.Range("T438").Value = 15
